# Contest 33 DC vs SRH.
# Fill in the 9 player scores for row 45 (Contest 33, "DC vs SRH").
# The points columns (D, G, J, M, P, S, V, Y, AB) already contain formulas
# that compute ranking-based points from these score cells, so updating the
# score cells alone is sufficient for the dependent formulas to recalculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E45").Value = 70
$ws.Range("H45").Value = 50
$ws.Range("K45").Value = 0
$ws.Range("N45").Value = 30
$ws.Range("Q45").Value = 80
$ws.Range("T45").Value = 40
$ws.Range("W45").Value = 20
$ws.Range("Z45").Value = 60
$ws.Range("AC45").Value = 100

$excel.CalculateFullRebuild()
